# Applies: updated "PERIOD TO EXPIRE" / "LAST UPDATE" figures (refreshed 8
# days later), a clearer "date is valid" remark on the Exam Dashboard, a
# wider COMMENTS column, and a unified bold/white header+title font across
# both dashboards.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------
# 1) Training Dashboard: "PERIOD TO EXPIRE" (H) shrinks by 8 days, and
#    "LAST UPDATE" (I) moves from 08-Sep-2025 to 16-Sep-2025.
# ---------------------------------------------------------------------
$periodUpdates = @{ 3 = 647; 4 = 435; 5 = 423; 6 = 677; 7 = 86; 8 = 182 }

foreach ($row in $periodUpdates.Keys) {
    $ws1.Range("H$row").Value = $periodUpdates[$row]

    # Assign as a text formula, then flatten to a static value via
    # copy/paste-values -- this keeps the cell a plain string (matching the
    # existing "DD-Mon-YYYY" text cells) instead of letting Excel's
    # autoconvert turn it into a real date serial number.
    $cell = $ws1.Range("I$row")
    $cell.Formula = "=""16-Sep-2025"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Exam Dashboard: widen the COMMENTS column and reword the remarks.
# ---------------------------------------------------------------------
$ws2.Range("E1").EntireColumn.ColumnWidth = 15 - (5/6)

$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"

# ---------------------------------------------------------------------
# 3) Restyle the title + header bands on both sheets: bold white text
#    (title drops its old 14pt size down to the regular 11pt).
# ---------------------------------------------------------------------
foreach ($sheetInfo in @(
        @{ Sheet = $ws1; HeaderRange = "A2:K2" },
        @{ Sheet = $ws2; HeaderRange = "A2:G2" }
    )) {
    $sheet = $sheetInfo.Sheet

    $title = $sheet.Range("A1")
    $title.Font.Bold = $true
    $title.Font.Size = 11
    $title.Font.Color = 16777215

    $header = $sheet.Range($sheetInfo.HeaderRange)
    $header.Font.Bold = $true
    $header.Font.Color = 16777215
}
